$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '25.894.20'
Set-TextValue $ws.Cells.Item(2, 5) '  -2.22%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '1.754.10'
Set-TextValue $ws.Cells.Item(3, 5) '  -4.58%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '1.000'
Set-TextValue $ws.Cells.Item(4, 5) '  -0.03%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '239.60'
Set-TextValue $ws.Cells.Item(5, 5) '  -8.00%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 5) '  -0.06%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.5100'
Set-TextValue $ws.Cells.Item(7, 5) '  -5.18%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '42.27'
Set-TextValue $ws.Cells.Item(8, 5) '  -5.77%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.2764'
Set-TextValue $ws.Cells.Item(9, 5) '  -5.27%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.06210'
Set-TextValue $ws.Cells.Item(10, 5) '  -10.46%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '1.748.18'
Set-TextValue $ws.Cells.Item(11, 5) '  -4.97%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 2) 'Solana'
Set-TextValue $ws.Cells.Item(12, 3) 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Cells.Item(12, 4) '15.81'
Set-TextValue $ws.Cells.Item(12, 5) '  -8.57%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 2) 'TRON'
Set-TextValue $ws.Cells.Item(13, 3) 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Cells.Item(13, 4) '0.06979'
Set-TextValue $ws.Cells.Item(13, 5) '  -2.99%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '0.6142'
Set-TextValue $ws.Cells.Item(14, 5) '  -15.43%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '4.537'
Set-TextValue $ws.Cells.Item(15, 5) '  -8.81%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '77.59'
Set-TextValue $ws.Cells.Item(16, 5) '  -12.81%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '0.9996'
Set-TextValue $ws.Cells.Item(17, 5) '  -0.13%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 5) '  -0.04%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '25.903.42'
Set-TextValue $ws.Cells.Item(19, 5) '  -2.26%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '0.000006932'
Set-TextValue $ws.Cells.Item(20, 5) '  -12.00%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '11.69'
Set-TextValue $ws.Cells.Item(21, 5) '  -15.13%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '1.970.67'
Set-TextValue $ws.Cells.Item(22, 5) '  -5.28%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '4.089'
Set-TextValue $ws.Cells.Item(23, 5) '  -10.75%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '5.277'
Set-TextValue $ws.Cells.Item(24, 5) '  -11.79%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '8.245'
Set-TextValue $ws.Cells.Item(25, 5) '  -10.12%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '138.14'
Set-TextValue $ws.Cells.Item(26, 5) '  -2.51%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '1.491'
Set-TextValue $ws.Cells.Item(27, 5) '  -12.55%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 2) 'LidoDAOToken'
Set-TextValue $ws.Cells.Item(28, 3) 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Cells.Item(28, 4) '1.827'
Set-TextValue $ws.Cells.Item(28, 5) '  -15.62%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 2) 'EthereumClassic'
Set-TextValue $ws.Cells.Item(29, 3) 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(29, 4) '15.09'
Set-TextValue $ws.Cells.Item(29, 5) '  -10.71%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '103.75'
Set-TextValue $ws.Cells.Item(30, 5) '  -6.46%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '0.08224'
Set-TextValue $ws.Cells.Item(31, 5) '  -7.28%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '3.711'
Set-TextValue $ws.Cells.Item(32, 5) '  -12.19%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '3.500'
Set-TextValue $ws.Cells.Item(33, 5) '  -12.92%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '0.04547'
Set-TextValue $ws.Cells.Item(34, 5) '  -5.94%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 5) '  -0.02%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '2.643'
Set-TextValue $ws.Cells.Item(36, 5) '  -9.19%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '0.9952'
Set-TextValue $ws.Cells.Item(37, 5) '  -11.78%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.6140'
Set-TextValue $ws.Cells.Item(38, 5) '  -15.01%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '2.719'
Set-TextValue $ws.Cells.Item(39, 5) '  -12.13%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '0.01559'
Set-TextValue $ws.Cells.Item(40, 5) '  -8.80%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '104.12'
Set-TextValue $ws.Cells.Item(41, 5) '  -2.57%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '1.000'
Set-TextValue $ws.Cells.Item(42, 5) '  +0.00%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '1.893'
Set-TextValue $ws.Cells.Item(43, 5) '  -17.52%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.3885'
Set-TextValue $ws.Cells.Item(44, 5) '  -16.66%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '0.7429'
Set-TextValue $ws.Cells.Item(45, 5) '  -17.64%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '4.948'
Set-TextValue $ws.Cells.Item(46, 5) '  -15.55%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '0.05426'
Set-TextValue $ws.Cells.Item(47, 5) '  -5.62%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '0.1117'
Set-TextValue $ws.Cells.Item(48, 5) '  -10.00%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) '6.030'
Set-TextValue $ws.Cells.Item(49, 5) '  -18.33%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '30.16'
Set-TextValue $ws.Cells.Item(50, 5) '  -13.20%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '52.91'
Set-TextValue $ws.Cells.Item(51, 5) '  -11.63%  '
